# Update the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions refresh). For numeric-looking price strings we
# prefix with an apostrophe so Excel keeps them as text (preserving exact
# formatting such as trailing zeros / thousands separators), then clear
# the resulting quote-prefix formatting so the cell style stays default.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.022.37"
$ws.Range("E2").Value = "  +6.56%  "
$ws.Range("D3").Value = "3.364.34"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'414.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.90%  "
$ws.Range("D6").Value = "'111.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "'0.589"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.645"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").Value = "'39.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").Value = "'0.0998"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "3.899.33"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "'20.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.72%  "
$ws.Range("D15").Value = "'8.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "3.352.94"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "60.876.78"
$ws.Range("E18").Value = "  +6.70%  "
$ws.Range("D19").Value = "'10.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'3.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +6.06%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'303.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'75.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("D26").Value = "'28.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").Value = "'7.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.33%  "
$ws.Range("D28").Value = "'4.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.181"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.62%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'7.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  +4.95%  "
$ws.Range("D32").Value = "'2.62"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +23.33%  "
$ws.Range("D33").Value = "'11.48"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.25%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'39.57"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").Value = "'52.46"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E41").Value = "  +7.05%  "
$ws.Range("D42").Value = "'137.54"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'3.97"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "'16.91"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").Value = "'22.64"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("E48").Value = "  +8.50%  "
$ws.Range("D49").Value = "2.191.16"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "'2.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "'1.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.25%  "
